$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventario")

# Update "valor_unitario" (column D) values for rows 3-6
$ws.Range("D3").Value = "500"
$ws.Range("D4").Value = "600"
$ws.Range("D5").Value = "200"
$ws.Range("D6").Value = "300"

# Update the view: zoom level and selected range
$ws.Activate()
$excel.ActiveWindow.Zoom = 140
$ws.Range("F3:F4").Select()
